# Locate the "(skype sarmkadans)" text and extend it with a date stamp,
# splitting it into the same run layout Word itself would produce
# (a standalone run for the flagged word "sarmkadans", plus a new
# trailing run holding the appended " 23.04.16y").

$d = $word.ActiveDocument

$findRng = $d.Content
$findRng.Find.Execute(" (skype sarmkadans)", $false, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)

$start = $findRng.Start

# Rewrite the matched text with the date suffix appended.
$findRng.Text = " (skype sarmkadans) 23.04.16y"

# Re-derive the offsets of "sarmkadans" and the appended " 23.04.16y"
# relative to the (stable) start of the range.
$wordStart = $start + 8
$wordEnd = $wordStart + 10          # "sarmkadans" is 10 characters
$suffixStart = $wordEnd + 1         # skip the ")" 
$suffixEnd = $suffixStart + 10      # " 23.04.16y" is 10 characters

# Force a run split around "sarmkadans" (becomes its own run, matching
# the rest of the formatting) by toggling Bold on and back off again.
$wordSub = $d.Range($wordStart, $wordEnd)
$wordSub.Bold = 1
$wordSub.Bold = 0

# Force a run split around the newly appended " 23.04.16y" the same way.
$suffixSub = $d.Range($suffixStart, $suffixEnd)
$suffixSub.Bold = 1
$suffixSub.Bold = 0
